$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 376; existing rows 376-474 shift down to 377-475.
$ws.Rows.Item(376).Insert()

# Populate the newly inserted row 376 with the new weekly price record.
$ws.Cells.Item(376, 1).Value  = 4
$ws.Cells.Item(376, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(376, 3).Value  = "Los Lagos"
$ws.Cells.Item(376, 4).Value  = 45135
$ws.Cells.Item(376, 5).Value  = 10
$ws.Cells.Item(376, 6).Value  = 100112043
$ws.Cells.Item(376, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(376, 8).Value  = "Sin especificar"
$ws.Cells.Item(376, 9).Value  = "Primera"
$ws.Cells.Item(376, 10).Value = 600
$ws.Cells.Item(376, 11).Value = 15000
$ws.Cells.Item(376, 12).Value = 15000
$ws.Cells.Item(376, 13).Value = 15000
$ws.Cells.Item(376, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(376, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(376, 16).Value = 250
$ws.Cells.Item(376, 17).Value = 60
$ws.Cells.Item(376, 18).Value = "Hortaliza"
